$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 "All Orders": a brand-new order (#9) came in at 2026-01-13 19:05.
# It belongs at the top of the data (row 2), pushing every existing order
# down by one row.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Insert a blank row at row 2 - this shifts the previous rows 2..9 down to
# rows 3..10 (their cell contents/types travel with them automatically).
$ws.Rows.Item(2).Insert()

# Duplicate the row that is now directly below (the old row 2, now row 3)
# into the new row 2 so the new row starts out with the same cell
# types/formatting (in particular so text-like values such as the phone
# number keep being stored as text instead of Excel auto-detecting a number).
$ws.Range("A3:L3").Copy($ws.Range("A2:L2"))

# Now overwrite the cells that actually hold new data for order #9.
$ws.Cells.Item(2,1).Value = 9
$ws.Cells.Item(2,2).Value = "2026-01-13 19:05"
$ws.Cells.Item(2,3).Value = "Sagar Borse"
# Column D (Phone) already holds the right text value "7588930329" thanks to
# the template copy above, so it is intentionally left untouched here -
# re-assigning it would make Excel re-detect it as a number.
$ws.Cells.Item(2,5).Value = "Test,"
$ws.Cells.Item(2,6).Value = "Stainless Steel Grater x1, Square Heat Pad x1, Circle Heat Pad x1, Leaf Tray (Large) x1"
$ws.Cells.Item(2,7).Value = 185
$ws.Cells.Item(2,8).Value = "NEW"
# Columns I (Payment), J (Notes), K (Cancel Reason) and L (Feedback) keep the
# values copied from the template row above (PENDING / blank / blank / blank).

# ---------------------------------------------------------------------------
# Sheet 2 "Daily Summary": roll the new order into the 2026-01-13 totals.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,2).Value = 9     # Total Orders: 8 -> 9
$ws2.Cells.Item(2,5).Value = 560   # Revenue: 375 -> 560
$ws2.Cells.Item(2,7).Value = 560   # Pending: 375 -> 560
